# Q3 Update - 2025
# Refresh the UNHCR Belize ("UN-BZE") dataset:
#  - drop the stray "Belgium" row (country of origin) for 2024
#  - drop the trailing "Unknown" row for 2024
#  - refresh the short-url stamp for every data row
#  - renumber the sequential "items" id column
#  - update refugees/asylum_seekers figures for the remaining 2024 rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the 2024 "Belgium" row (worksheet row 184). Everything below
#    shifts up by one row automatically.
$ws.Rows(184).Delete()

# 2) Remove the trailing 2024 "Unknown" row. After the delete above, the
#    sheet now ends at row 190 (was 191), so the "Unknown" row now sits at
#    row 190.
$ws.Rows(190).Delete()

# 3) The page's short-url stamp changed for this refresh; it is repeated
#    down column B for every data row (rows 2-189 after the two deletes).
$ws.Range("B2:B189").Value2 = "Avz9E1"

# 4) Renumber the sequential "items" id (column D) for the six remaining
#    2024 rows so the ids stay contiguous (183-188).
$ws.Range("D184").Value2 = "183"
$ws.Range("D185").Value2 = "184"
$ws.Range("D186").Value2 = "185"
$ws.Range("D187").Value2 = "186"
$ws.Range("D188").Value2 = "187"
$ws.Range("D189").Value2 = "188"

# 5) Apply the refreshed refugees (N) / asylum_seekers (O) figures for the
#    remaining 2024 rows (Cuba, Guatemala, Haiti, Honduras, Nicaragua, El
#    Salvador).
$ws.Range("N184").Value2 = 0
$ws.Range("O184").Value2 = 5

$ws.Range("N185").Value2 = 19
$ws.Range("O185").Value2 = 418

$ws.Range("N186").Value2 = 0
$ws.Range("O186").Value2 = 9

$ws.Range("N187").Value2 = 95
$ws.Range("O187").Value2 = 438

$ws.Range("N188").Value2 = 16
$ws.Range("O188").Value2 = 6

$ws.Range("N189").Value2 = 156
$ws.Range("O189").Value2 = 1125
